$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range so stale cells (e.g. old column E data,
# old row 3/4/5 values) don't linger when the new layout is written.
$ws.Range("A1:I7").Clear()

# --- Header row ---
# Shared-string table order matters: new strings are appended to
# xl/sharedStrings.xml in first-write order. The target file has
# "波长（nm）" at index 6 and "是否为有源区（是为1）" at index 7, so
# I1 must be written before D1.
$ws.Range("A1").Value = "层厚（nm）"
$ws.Range("B1").Value = "折射率实部"
$ws.Range("C1").Value = "吸收（/cm）"
$ws.Range("F1").Value = "波导折射率实部"
$ws.Range("G1").Value = "cladding折射率实部"
$ws.Range("H1").Value = "波导厚度（nm）"
$ws.Range("I1").Value = "波长（nm）"
$ws.Range("D1").Value = "是否为有源区（是为1）"

# --- Row 2 ---
$ws.Range("A2").Value = 500
$ws.Range("B2").Value = 2.4500000000000002
$ws.Range("C2").Value = 50
$ws.Range("F2").Value = 2.5
$ws.Range("G2").Value = 2.4500000000000002
$ws.Range("H2").Value = 600
$ws.Range("I2").Value = 450

# --- Row 3 ---
$ws.Range("A3").Value = 300
$ws.Range("B3").Value = 2.5
$ws.Range("C3").Value = 0

# --- Row 4 ---
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 2.6
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

# --- Row 5 ---
$ws.Range("A5").Value = 300
$ws.Range("B5").Value = 2.5
$ws.Range("C5").Value = 0

# --- Row 6 ---
$ws.Range("A6").Value = 500
$ws.Range("B6").Value = 2.4500000000000002
$ws.Range("C6").Value = 0

# --- Row 7 ---
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0

# Selection matches the target diff (active cell D4)
$ws.Range("D4").Select()
